$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking strings need to be
# forced to text (via NumberFormat "@") so Excel keeps them as strings
# instead of auto-converting to numbers.
$textCells = @(
    "D5",
    "D6",
    "D9",
    "D10",
    "D12",
    "D13",
    "D15",
    "D17",
    "D19",
    "D21",
    "D23",
    "D24",
    "D25",
    "D27",
    "D29",
    "D30",
    "D31",
    "D32",
    "D34",
    "D35",
    "D36",
    "D37",
    "D39",
    "D40",
    "D41",
    "D43",
    "D44",
    "D47",
    "D48",
    "D50",
    "D51"
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values (numeric-looking, now text-formatted cells)
$ws.Range("D5").Value2 = '352.13'
$ws.Range("D6").Value2 = '105.26'
$ws.Range("D9").Value2 = '0.595'
$ws.Range("D10").Value2 = '37.35'
$ws.Range("D12").Value2 = '0.0845'
$ws.Range("D13").Value2 = '18.74'
$ws.Range("D15").Value2 = '7.41'
$ws.Range("D17").Value2 = '0.976'
$ws.Range("D19").Value2 = '3.31'
$ws.Range("D21").Value2 = '13.21'
$ws.Range("D23").Value2 = '68.70'
$ws.Range("D24").Value2 = '264.26'
$ws.Range("D25").Value2 = '2.67'
$ws.Range("D27").Value2 = '26.42'
$ws.Range("D29").Value2 = '7.18'
$ws.Range("D30").Value2 = '0.108'
$ws.Range("D31").Value2 = '6.22'
$ws.Range("D32").Value2 = '10.01'
$ws.Range("D34").Value2 = '35.31'
$ws.Range("D35").Value2 = '50.73'
$ws.Range("D36").Value2 = '0.998'
$ws.Range("D37").Value2 = '0.0425'
$ws.Range("D39").Value2 = '2.82'
$ws.Range("D40").Value2 = '17.21'
$ws.Range("D41").Value2 = '1.90'
$ws.Range("D43").Value2 = '22.93'
$ws.Range("D44").Value2 = '120.98'
$ws.Range("D47").Value2 = '3.23'
$ws.Range("D48").Value2 = '2.31'
$ws.Range("D50").Value2 = '0.237'
$ws.Range("D51").Value2 = '0.0319'

# Apply the updated values for cells that remain safely text as-is
$ws.Range("D2").Value2 = '51.723.45'
$ws.Range("E2").Value2 = '  -0.44%  '
$ws.Range("D3").Value2 = '2.936.64'
$ws.Range("E3").Value2 = '  +0.25%  '
$ws.Range("E4").Value2 = '  +0.01%  '
$ws.Range("E5").Value2 = '  -1.64%  '
$ws.Range("E6").Value2 = '  -4.73%  '
$ws.Range("E7").Value2 = '  -3.80%  '
$ws.Range("E9").Value2 = '  -6.08%  '
$ws.Range("E10").Value2 = '  -5.30%  '
$ws.Range("E11").Value2 = '  +1.82%  '
$ws.Range("E12").Value2 = '  -3.92%  '
$ws.Range("E13").Value2 = '  -4.87%  '
$ws.Range("D14").Value2 = '3.399.39'
$ws.Range("E15").Value2 = '  -6.87%  '
$ws.Range("D16").Value2 = '2.928.80'
$ws.Range("E16").Value2 = '  -0.02%  '
$ws.Range("E17").Value2 = '  -0.98%  '
$ws.Range("D18").Value2 = '51.602.71'
$ws.Range("E18").Value2 = '  -0.71%  '
$ws.Range("E19").Value2 = '  -2.07%  '
$ws.Range("E20").Value2 = '  -4.22%  '
$ws.Range("E21").Value2 = '  -6.34%  '
$ws.Range("D22").Value2 = '0.0₃0953'
$ws.Range("E22").Value2 = '  -3.15%  '
$ws.Range("E23").Value2 = '  -3.39%  '
$ws.Range("E24").Value2 = '  -2.26%  '
$ws.Range("E25").Value2 = '  -5.85%  '
$ws.Range("E26").Value2 = '  -6.55%  '
$ws.Range("E27").Value2 = '  -2.70%  '
$ws.Range("E28").Value2 = '  +0.07%  '
$ws.Range("E29").Value2 = '  -4.87%  '
$ws.Range("E30").Value2 = '  +0.41%  '
$ws.Range("E31").Value2 = '  +2.19%  '
$ws.Range("E32").Value2 = '  -5.80%  '
$ws.Range("E33").Value2 = '  -5.50%  '
$ws.Range("E34").Value2 = '  -7.53%  '
$ws.Range("E35").Value2 = '  -2.93%  '
$ws.Range("E36").Value2 = '  -0.09%  '
$ws.Range("E37").Value2 = '  -4.31%  '
$ws.Range("E38").Value2 = '  -1.60%  '
$ws.Range("E39").Value2 = '  +2.62%  '
$ws.Range("E40").Value2 = '  -6.89%  '
$ws.Range("E41").Value2 = '  -5.89%  '
$ws.Range("E42").Value2 = '  -4.50%  '
$ws.Range("E43").Value2 = '  -1.49%  '
$ws.Range("E44").Value2 = '  +1.30%  '
$ws.Range("E45").Value2 = '  -0.34%  '
$ws.Range("D46").Value2 = '2.096.85'
$ws.Range("E46").Value2 = '  -1.93%  '
$ws.Range("B47").Value2 = 'NEARProtocol'
$ws.Range("C47").Value2 = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("E47").Value2 = '  -7.30%  '
$ws.Range("B48").Value2 = 'ApeXProtocol'
$ws.Range("C48").Value2 = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("E48").Value2 = '  -7.17%  '
$ws.Range("D49").Value2 = '3.232.31'
$ws.Range("E49").Value2 = '  +0.38%  '
$ws.Range("E50").Value2 = '  -4.90%  '
$ws.Range("E51").Value2 = '  -5.00%  '
